# Re-sync match-odds rows that were reordered upstream (league base update 12-06-2024).
# Each group of row numbers below had their full data payload (everything except the
# running index in column A, the competition name in column C and the match Date in
# column D) rotated among the listed rows, in the order given.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Group 1: rows 20, 21
$rowData = @{}
$rowData[20] = @{ "B" = 6920350; "E" = "Phonix Lubeck"; "F" = "SSV Jeddeloh"; "G" = 7; "H" = 0; "I" = 3; "J" = 0; "K" = "H"; "L" = 2.25; "M" = 3.75; "N" = 2.5; "O" = 2.05; "P" = 4.2; "Q" = 2.75; "R" = -0.25; "S" = 1.875; "T" = 1.925; "U" = 3; "V" = 1.875; "W" = 1.925; "X" = 1.05; "Y" = -1; "Z" = -1; "AA" = 0.875; "AB" = -1; "AC" = 0.875; "AD" = -1 }
$rowData[21] = @{ "B" = 6920351; "E" = "Hamburg SV II"; "F" = "SV DrochtersenAssel"; "G" = 1; "H" = 1; "I" = 1; "J" = 0; "K" = "D"; "L" = 3; "M" = 3.6; "N" = 2; "O" = 5; "P" = 4; "Q" = 1.615; "R" = 1; "S" = 1.775; "T" = 2.025; "U" = 2.75; "V" = 1.775; "W" = 2.025; "X" = -1; "Y" = 3; "Z" = -1; "AA" = 0.7749999999999999; "AB" = -1; "AC" = -1; "AD" = 1.025 }
$rows = @(20, 21)
$sourceForTarget = @{ 20 = 21; 21 = 20 }
foreach ($targetRow in $rows) {
    $srcRow = $sourceForTarget[$targetRow]
    $src = $rowData[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }
}

# Group 2: rows 35, 36
$rowData = @{}
$rowData[35] = @{ "B" = 6920365; "E" = "Hamburg SV II"; "F" = "Bremer SV"; "G" = 3; "H" = 1; "I" = 1; "J" = 1; "K" = "H"; "L" = 2.55; "M" = 3.75; "N" = 2.2; "O" = 2.1; "P" = 3.6; "Q" = 2.75; "R" = -0.25; "S" = 1.925; "T" = 1.875; "U" = 2.5; "V" = 1.825; "W" = 1.975; "X" = 1.1; "Y" = -1; "Z" = -1; "AA" = 0.925; "AB" = -1; "AC" = 0.825; "AD" = -1 }
$rowData[36] = @{ "B" = 6919326; "E" = "FC Teutonia 05"; "F" = "SV DrochtersenAssel"; "G" = 6; "H" = 0; "I" = 4; "J" = 0; "K" = "H"; "L" = 1.5; "M" = 3.75; "N" = 5.5; "O" = 1.85; "P" = 3.3; "Q" = 3.5; "R" = -0.5; "S" = 1.925; "T" = 1.875; "U" = 3; "V" = 1.9; "W" = 1.9; "X" = 0.8500000000000001; "Y" = -1; "Z" = -1; "AA" = 0.925; "AB" = -1; "AC" = 0.8999999999999999; "AD" = -1 }
$rows = @(35, 36)
$sourceForTarget = @{ 35 = 36; 36 = 35 }
foreach ($targetRow in $rows) {
    $srcRow = $sourceForTarget[$targetRow]
    $src = $rowData[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }
}

# Group 3: rows 95, 97, 96
$rowData = @{}
$rowData[95] = @{ "B" = 6920417; "E" = "SV DrochtersenAssel"; "F" = "SSV Jeddeloh"; "G" = 2; "H" = 2; "I" = 1; "J" = 0; "K" = "D"; "L" = 2.2; "M" = 3.75; "N" = 2.55; "O" = 1.727; "P" = 4; "Q" = 3.6; "R" = -0.75; "S" = 2.025; "T" = 1.825; "U" = 2.75; "V" = 1.85; "W" = 2; "X" = -1; "Y" = 3; "Z" = -1; "AA" = -1; "AB" = 0.825; "AC" = 0.8500000000000001; "AD" = -1 }
$rowData[97] = @{ "B" = 6920415; "E" = "SC SpelleVenhaus"; "F" = "TSV Havelse"; "G" = 1; "H" = 3; "I" = 1; "J" = 1; "K" = "A"; "L" = 2.5; "M" = 3.75; "N" = 2.25; "O" = 4.5; "P" = 4.2; "Q" = 1.5; "R" = 1; "S" = 1.95; "T" = 1.9; "U" = 3.25; "V" = 2; "W" = 1.85; "X" = -1; "Y" = -1; "Z" = 0.5; "AA" = -1; "AB" = 0.8999999999999999; "AC" = 1; "AD" = -1 }
$rowData[96] = @{ "B" = 6920416; "E" = "Bremer SV"; "F" = "VfB Oldenburg"; "G" = 0; "H" = 2; "I" = 0; "J" = 0; "K" = "A"; "L" = 2.15; "M" = 3.75; "N" = 2.625; "O" = 4.5; "P" = 4; "Q" = 1.55; "R" = 1; "S" = 1.8; "T" = 2.05; "U" = 3; "V" = 1.975; "W" = 1.875; "X" = -1; "Y" = -1; "Z" = 0.55; "AA" = -1; "AB" = 1.05; "AC" = -1; "AD" = 0.875 }
$rows = @(95, 97, 96)
$sourceForTarget = @{ 95 = 97; 97 = 96; 96 = 95 }
foreach ($targetRow in $rows) {
    $srcRow = $sourceForTarget[$targetRow]
    $src = $rowData[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }
}

# Group 4: rows 121, 122, 123
$rowData = @{}
$rowData[121] = @{ "B" = 6920445; "E" = "VfB Oldenburg"; "F" = "Eimsbutteler TV"; "G" = 4; "H" = 1; "I" = 2; "J" = 1; "K" = "H"; "L" = 1.363; "M" = 4.75; "N" = 6; "O" = 1.2; "P" = 5.75; "Q" = 9; "R" = -2; "S" = 1.95; "T" = 1.9; "U" = 3.75; "V" = 1.975; "W" = 1.875; "X" = 0.2; "Y" = -1; "Z" = -1; "AA" = 0.95; "AB" = -1; "AC" = 0.9750000000000001; "AD" = -1 }
$rowData[122] = @{ "B" = 6920446; "E" = "Bremer SV"; "F" = "St Pauli II"; "G" = 1; "H" = 1; "I" = 1; "J" = 1; "K" = "D"; "L" = 3; "M" = 4; "N" = 1.909; "O" = 3.1; "P" = 3.75; "Q" = 1.95; "R" = 0.5; "S" = 1.875; "T" = 1.975; "U" = 3.25; "V" = 1.95; "W" = 1.9; "X" = -1; "Y" = 2.75; "Z" = -1; "AA" = 0.875; "AB" = -1; "AC" = -1; "AD" = 0.8999999999999999 }
$rowData[123] = @{ "B" = 6920447; "E" = "Holstein Kiel II"; "F" = "TuS BlauWeiss Lohne"; "G" = 1; "H" = 1; "I" = 0; "J" = 0; "K" = "D"; "L" = 1.5; "M" = 4.2; "N" = 4.75; "O" = 1.6; "P" = 4; "Q" = 4.2; "R" = -1; "S" = 1.975; "T" = 1.825; "U" = 3.5; "V" = 1.925; "W" = 1.875; "X" = -1; "Y" = 3; "Z" = -1; "AA" = -1; "AB" = 0.825; "AC" = -1; "AD" = 0.875 }
$rows = @(121, 122, 123)
$sourceForTarget = @{ 121 = 122; 122 = 123; 123 = 121 }
foreach ($targetRow in $rows) {
    $srcRow = $sourceForTarget[$targetRow]
    $src = $rowData[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $src[$c]
    }
}

Write-Host "done"
